$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test case name in A4 from "Teste3" to "Invalid Search"
$ws.Range("A4").Value = "Invalid Search"

# Move the active selection to A3 (was C7)
$ws.Range("A3").Select()
